{"js": "// Word JS API (Office.js) edit script.\n// Body is the async (context) => { ... } function.\n//\n// Content edits reproduced from the diff:\n//  1. Heading \"Prezentarea sistemului\" -> \"System description\"\n//  2. \"The project aims to predict whether or not a COVID test ...\"\n//       -> \"The project aims to predict whether a COVID test ...\"   (drop \"or not\")\n//  3. \"Recall, precision , f1_score is for patients who are negative.\"\n//       -> \"Recall, precision, f1_score are for predictions on patients who are negative.\"\n//  4. Explicitly set the section's page orientation to Portrait\n//     (diff adds an explicit w:orient=\"portrait\" to <w:pgSz/>).\n\nconst body = context.document.body;\n\n// 1. Heading text replacement.\nconst headingHits = body.search(\"Prezentarea sistemului\", { matchCase: true });\nheadingHits.load(\"items\");\nawait context.sync();\nif (headingHits.items.length > 0) {\n  headingHits.items[0].insertText(\"System description\", \"Replace\");\n}\n\n// 2. Remove \"or not\" from the project-aims sentence.\nconst orNotHits = body.search(\"whether or not a COVID test\", { matchCase: true });\norNotHits.load(\"items\");\nawait context.sync();\nif (orNotHits.items.length > 0) {\n  orNotHits.items[0].insertText(\"whether a COVID test\", \"Replace\");\n}\n\n// 3. Update the recall/precision/f1_score sentence.\nconst metricsHits = body.search(\n  \"Recall, precision , f1_score is for patients who are negative.\",\n  { matchCase: true }\n);\nmetricsHits.load(\"items\");\nawait context.sync();\nif (metricsHits.items.length > 0) {\n  metricsHits.items[0].insertText(\n    \"Recall, precision, f1_score are for predictions on patients who are negative.\",\n    \"Replace\"\n  );\n}\n\n// 4. Make the (already-portrait) page orientation explicit.\nconst section = context.document.sections.getFirst();\nsection.pageSetup.orientation = \"Portrait\";\n\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# $word.ActiveDocument is the open document.\n#\n# Content edits reproduced from the diff:\n#  1. Heading \"Prezentarea sistemului\" -> \"System description\"\n#  2. \"The project aims to predict whether or not a COVID test ...\"\n#       -> \"The project aims to predict whether a COVID test ...\"   (drop \"or not\")\n#  3. \"Recall, precision , f1_score is for patients who are negative.\"\n#       -> \"Recall, precision, f1_score are for predictions on patients who are negative.\"\n#  4. Explicitly set the page orientation to Portrait\n#     (diff adds an explicit w:orient=\"portrait\" to <w:pgSz/>).\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# 1. Heading text replacement.\nReplace-Text \"Prezentarea sistemului\" \"System description\"\n\n# 2. Remove \"or not\" from the project-aims sentence.\nReplace-Text \"whether or not a COVID test\" \"whether a COVID test\"\n\n# 3. Update the recall/precision/f1_score sentence.\nReplace-Text \"Recall, precision , f1_score is for patients who are negative.\" \"Recall, precision, f1_score are for predictions on patients who are negative.\"\n\n# 4. Make the (already-portrait) page orientation explicit.\n# wdOrientPortrait = 0\n$d.PageSetup.Orientation = 0\n"}
